$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.230310916900635
$ws.Range("B1").Value = 2.172172784805298
$ws.Range("C1").Value = 4.168947219848633
$ws.Range("D1").Value = 3.072719097137451
$ws.Range("E1").Value = 1.075792789459229
